# Update quizvragen via Admin
#
# - Sheet "DC": remove the stray/incomplete row 3 ("Wat betekend
#   gelijkstroom") so the sheet only keeps the header + the q3 (vermogen)
#   question. Dimension collapses from A1:L3 to A1:L2.
# - Sheet "Wiskunde 3": remove row 3 (the old q2 "cos(α)" question). All
#   following rows (old q3..q20) shift up one row, which also naturally
#   renumbers their ids down by one (old q3 -> new q3 occupying row 3,
#   old q4 -> new q4 on row 4, etc). Dimension collapses from A1:L21 to
#   A1:L20.
# - Sheet "AC": untouched.

$wb = $excel.ActiveWorkbook

$wsDC = $wb.Worksheets.Item("DC")
$wsDC.Rows(3).Delete()

$wsWiskunde = $wb.Worksheets.Item("Wiskunde 3")
$wsWiskunde.Rows(3).Delete()
